# Fix i/o errors: rename "Source Name" columns to "Sample Name"
# on the "Events-Harvest" sheet/table (annotationTable).
#
# The table's header row cells hold the column names, so updating the
# cell values keeps both the worksheet cells and the underlying Excel
# Table (ListObject) column names in sync.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Events-Harvest")

# "Input [Source Name]" -> "Input [Sample Name]" (first table column, cell A1)
$ws.Range("A1").Value2 = "Input [Sample Name]"

# "Output [Source Name]" -> "Output [Sample Name]" (last table column, cell AF1)
$ws.Range("AF1").Value2 = "Output [Sample Name]"
